$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "s + + heet1"

# Update the header text in C1 (shared string "Total Case" -> "Total Case to date 13/04/2023")
$ws.Range("C1").Value = "Total Case to date 13/04/2023"

# Update column C values (rows 2-22)
$ws.Range("C2").Value = 4143150
$ws.Range("C3").Value = 2712436
$ws.Range("C4").Value = 2456161
$ws.Range("C5").Value = 2402802
$ws.Range("C6").Value = 2146082
$ws.Range("C7").Value = 1822963
$ws.Range("C8").Value = 1728126
$ws.Range("C9").Value = 1632702
$ws.Range("C10").Value = 1597294
$ws.Range("C11").Value = 716098
$ws.Range("C12").Value = 664258
$ws.Range("C13").Value = 655155
$ws.Range("C14").Value = 634472
$ws.Range("C15").Value = 578199
$ws.Range("C16").Value = 512209
$ws.Range("C17").Value = 441114
$ws.Range("C18").Value = 295517
$ws.Range("C19").Value = 245446
$ws.Range("C20").Value = 200156
$ws.Range("C21").Value = 102165
$ws.Range("C22").Value = 50665
